# Weekly update: a new price-report row is inserted at the top of the data
# (row 4), pushing all the existing observations (old rows 4-40) down by one
# row (new rows 5-41). The new row 4 is a fresh entry for "Macroferia
# Regional de Talca" / "Arándano (blue)" dated 2021-12-07 (serial 44537).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 4; everything below (old rows 4-40) shifts
# down to rows 5-41, carrying its formatting/styles with it.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new observation.
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Macroferia Regional de Talca"
$ws.Range("C4").Value = "Maule"
$ws.Range("D4").Value = 44537
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101001
$ws.Range("J4").Value = "Arándano (blue)"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 3600
$ws.Range("O4").Value = 3600
$ws.Range("P4").Value = 3600
$ws.Range("Q4").Value = '$/bandeja 2 kilos'
$ws.Range("R4").Value = "Provincia de Linares"
$ws.Range("S4").Value = 1800
$ws.Range("T4").Value = 2
